# feat: add file read me
# Append new "question" rows (74-82) to the Reviews sheet, replicating the
# question/answer export rows that were added to exportQuestion.xlsx.
# Column layout: A=question text, B=teacher id, C=subject, D=difficulty,
# E..H=answer choices. A choice cell that represents the "correct answer"
# is highlighted with a yellow fill (same indexed color already used by the
# existing rows in the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ row = 74; cells = @(@{ col = "A"; val = "1 + 3 = ?"; isText = $true; highlight = $false }; @{ col = "B"; val = 16.0; isText = $false; highlight = $false }; @{ col = "C"; val = "Toán cao cấp"; isText = $true; highlight = $false }; @{ col = "D"; val = "EASY"; isText = $true; highlight = $false }; @{ col = "E"; val = "4.0"; isText = $true; highlight = $true }; @{ col = "F"; val = "ok"; isText = $true; highlight = $false }) },
    @{ row = 75; cells = @(@{ col = "A"; val = "Toi la ai ?"; isText = $true; highlight = $false }; @{ col = "B"; val = 16.0; isText = $false; highlight = $false }; @{ col = "C"; val = "Lập trình Java"; isText = $true; highlight = $false }; @{ col = "D"; val = "HARD"; isText = $true; highlight = $false }; @{ col = "E"; val = "loc"; isText = $true; highlight = $false }; @{ col = "F"; val = "huy"; isText = $true; highlight = $false }; @{ col = "G"; val = "hung"; isText = $true; highlight = $false }; @{ col = "H"; val = "khong co dap dan dung"; isText = $true; highlight = $true }) },
    @{ row = 76; cells = @(@{ col = "A"; val = "cau hoi 1 ne"; isText = $true; highlight = $false }; @{ col = "B"; val = 16.0; isText = $false; highlight = $false }; @{ col = "C"; val = "Toán cao cấp"; isText = $true; highlight = $false }; @{ col = "D"; val = "EASY"; isText = $true; highlight = $false }; @{ col = "E"; val = "fsf"; isText = $true; highlight = $false }; @{ col = "F"; val = "fa"; isText = $true; highlight = $true }; @{ col = "G"; val = "fsd"; isText = $true; highlight = $false }; @{ col = "H"; val = "fsd"; isText = $true; highlight = $false }) },
    @{ row = 77; cells = @(@{ col = "A"; val = "cau hoi 2"; isText = $true; highlight = $false }; @{ col = "B"; val = 16.0; isText = $false; highlight = $false }; @{ col = "C"; val = "Toán cao cấp"; isText = $true; highlight = $false }; @{ col = "D"; val = "EASY"; isText = $true; highlight = $false }; @{ col = "E"; val = "fád"; isText = $true; highlight = $false }; @{ col = "F"; val = "fsa"; isText = $true; highlight = $false }; @{ col = "G"; val = "fsda"; isText = $true; highlight = $true }; @{ col = "H"; val = "fsdfsf"; isText = $true; highlight = $false }) },
    @{ row = 78; cells = @(@{ col = "A"; val = "1 + 3 = ?"; isText = $true; highlight = $false }; @{ col = "B"; val = 16.0; isText = $false; highlight = $false }; @{ col = "C"; val = "Toán cao cấp"; isText = $true; highlight = $false }; @{ col = "D"; val = "EASY"; isText = $true; highlight = $false }; @{ col = "E"; val = "4.0"; isText = $true; highlight = $true }; @{ col = "F"; val = "ok"; isText = $true; highlight = $false }) },
    @{ row = 79; cells = @(@{ col = "A"; val = "Toi la ai ?"; isText = $true; highlight = $false }; @{ col = "B"; val = 16.0; isText = $false; highlight = $false }; @{ col = "C"; val = "Lập trình Java"; isText = $true; highlight = $false }; @{ col = "D"; val = "HARD"; isText = $true; highlight = $false }; @{ col = "E"; val = "loc"; isText = $true; highlight = $false }; @{ col = "F"; val = "huy"; isText = $true; highlight = $false }; @{ col = "G"; val = "hung"; isText = $true; highlight = $false }; @{ col = "H"; val = "khong co dap dan dung"; isText = $true; highlight = $true }) },
    @{ row = 80; cells = @(@{ col = "A"; val = "cau hoi so 78 = ?"; isText = $true; highlight = $false }; @{ col = "B"; val = 16.0; isText = $false; highlight = $false }; @{ col = "C"; val = "Toán cao cấp"; isText = $true; highlight = $false }; @{ col = "D"; val = "EASY"; isText = $true; highlight = $false }; @{ col = "E"; val = "32"; isText = $true; highlight = $true }; @{ col = "F"; val = "78"; isText = $true; highlight = $false }; @{ col = "G"; val = "213"; isText = $true; highlight = $false }) },
    @{ row = 81; cells = @(@{ col = "A"; val = "1 + 3 = ?"; isText = $true; highlight = $false }; @{ col = "B"; val = 16.0; isText = $false; highlight = $false }; @{ col = "C"; val = "Toán cao cấp"; isText = $true; highlight = $false }; @{ col = "D"; val = "EASY"; isText = $true; highlight = $false }; @{ col = "E"; val = "4.0"; isText = $true; highlight = $true }; @{ col = "F"; val = "ok"; isText = $true; highlight = $false }) },
    @{ row = 82; cells = @(@{ col = "A"; val = "Toi la ai ?"; isText = $true; highlight = $false }; @{ col = "B"; val = 16.0; isText = $false; highlight = $false }; @{ col = "C"; val = "Lập trình Java"; isText = $true; highlight = $false }; @{ col = "D"; val = "HARD"; isText = $true; highlight = $false }; @{ col = "E"; val = "loc"; isText = $true; highlight = $false }; @{ col = "F"; val = "huy"; isText = $true; highlight = $false }; @{ col = "G"; val = "hung"; isText = $true; highlight = $false }; @{ col = "H"; val = "khong co dap dan dung"; isText = $true; highlight = $true }) },
)

foreach ($r in $rows) {
    $rowNum = $r.row
    foreach ($cell in $r.cells) {
        $addr = "$($cell.col)$rowNum"
        $rng = $ws.Range($addr)

        if ($cell.isText -and ($cell.val -match '^-?\d+(\.\d+)?$')) {
            # Value looks numeric (e.g. "4.0", "78") but must be stored as
            # text, matching the answer-choice columns elsewhere in this
            # sheet (e.g. existing "3", "5", "8", "10" choices).
            $rng.NumberFormat = "@"
        }

        $rng.Value = $cell.val

        if ($cell.highlight) {
            # Same yellow fill (indexed color 13) used by the correct-answer
            # cells in the pre-existing rows of this sheet.
            $rng.Interior.Color = 65535
        }
    }
}
